$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.532.21"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3
$ws.Range("D3").Value = "1.826.72"
$ws.Range("E3").Value = "  -0.03%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "'317.07"
$ws.Range("E5").Value = "  +0.45%  "

# Row 6
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.5162"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "'0.3890"
$ws.Range("E8").Value = "  -1.30%  "

# Row 9
$ws.Range("D9").Value = "'0.08424"
$ws.Range("E9").Value = "  +8.92%  "

# Row 10
$ws.Range("D10").Value = "'1.122"
$ws.Range("E10").Value = "  +0.72%  "

# Row 11
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("D12").Value = "'6.433"
$ws.Range("E12").Value = "  +2.31%  "

# Row 13
$ws.Range("D13").Value = "'21.28"
$ws.Range("E13").Value = "  +0.85%  "

# Row 14
$ws.Range("E14").Value = "  +0.16%  "

# Row 15
$ws.Range("E15").Value = "  -0.63%  "

# Row 16
$ws.Range("D16").Value = "1.828.06"
$ws.Range("E16").Value = "  +0.25%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001135"
$ws.Range("E17").Value = "  +4.93%  "

# Row 18
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'94.35"
$ws.Range("E18").Value = "  +1.08%  "

# Row 19
$ws.Range("D19").Value = "'0.06625"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").Value = "'17.80"
$ws.Range("E20").Value = "  +0.58%  "

# Row 21
$ws.Range("E21").Value = "  +0.15%  "

# Row 22
$ws.Range("D22").Value = "'6.090"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23
$ws.Range("D23").Value = "28.578.16"
$ws.Range("E23").Value = "  +0.41%  "

# Row 24
$ws.Range("E24").Value = "  +2.63%  "

# Row 25
$ws.Range("E25").Value = "  +1.00%  "

# Row 26
$ws.Range("D26").Value = "'21.27"
$ws.Range("E26").Value = "  +3.04%  "

# Row 27
$ws.Range("D27").Value = "'159.84"
$ws.Range("E27").Value = "  +1.43%  "

# Row 28
$ws.Range("D28").Value = "2.036.34"
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").Value = "'2.409"
$ws.Range("E29").Value = "  -1.61%  "

# Row 30
$ws.Range("D30").Value = "'125.84"
$ws.Range("E30").Value = "  +0.68%  "

# Row 31
$ws.Range("D31").Value = "'0.1098"
$ws.Range("E31").Value = "  -0.16%  "

# Row 32
$ws.Range("D32").Value = "'1.100"
$ws.Range("E32").Value = "  -2.83%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.07668"
$ws.Range("E33").Value = "  +7.09%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.743"
$ws.Range("E34").Value = "  +1.39%  "

# Row 35
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("D36").Value = "'0.02395"
$ws.Range("E36").Value = "  +2.89%  "

# Row 37
$ws.Range("D37").Value = "'0.2231"
$ws.Range("E37").Value = "  -0.30%  "

# Row 38
$ws.Range("D38").Value = "'5.263"
$ws.Range("E38").Value = "  +1.98%  "

# Row 39
$ws.Range("D39").Value = "'8.747"
$ws.Range("E39").Value = "  -2.53%  "

# Row 40
$ws.Range("D40").Value = "'0.6379"
$ws.Range("E40").Value = "  +2.07%  "

# Row 41
$ws.Range("E41").Value = "  +1.82%  "

# Row 42
$ws.Range("D42").Value = "'1.189"
$ws.Range("E42").Value = "  -0.19%  "

# Row 43
$ws.Range("D43").Value = "'1.400"
$ws.Range("E43").Value = "  +0.41%  "

# Row 44
$ws.Range("D44").Value = "'13.58"
$ws.Range("E44").Value = "  +0.68%  "

# Row 45
$ws.Range("D45").Value = "'0.6054"
$ws.Range("E45").Value = "  +2.50%  "

# Row 46
$ws.Range("E46").Value = "  +2.05%  "

# Row 47
$ws.Range("D47").Value = "'127.36"
$ws.Range("E47").Value = "  +2.23%  "

# Row 48
$ws.Range("E48").Value = "  +0.98%  "

# Row 49
$ws.Range("E49").Value = "  +1.93%  "

# Row 51
$ws.Range("D51").Value = "'74.73"
$ws.Range("E51").Value = "  +1.14%  "
